# feat: add 2022-Q4 data
#
# The workbook has a "总计" (summary) sheet followed by one per-quarter
# sheet ("2021-Q4"). This adds a new "2022-Q4" quarter sheet (inserted
# right after "总计", before "2021-Q4") with its fund-holding data, and
# updates the summary sheet with a new row for 2022-Q4 (pushing the
# existing 2021-Q4 summary row down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet: insert the 2022-Q4 row above the
#    existing 2021-Q4 row (which moves from row 2 to row 3).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push the formatting of row 2 down into row 3 first, so the moved
# 2021-Q4 row keeps its original look (bordered/centered index cell).
$summary.Range("A2:D2").Copy()
$summary.Range("A3:D3").PasteSpecial(-4122)

# Existing 2021-Q4 data moves down to row 3; its index bumps 0 -> 1.
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2021-Q4"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.04

# New 2022-Q4 summary data goes into row 2.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.02

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计" (so it sits
#    before the existing "2021-Q4" sheet).
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$newSheet.Name = "2022-Q4"

# Match the sheet-level look & feel (outline defaults + page margins)
# used by the other sheets in this workbook.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Copy the header/index-column cell formatting from the 2021-Q4 sheet
# so the new sheet's styling matches the rest of the workbook. (Fetch
# this sheet reference only now, after the insert above, so it isn't
# left dangling against the pre-insert sheet collection.)
$quarter2021 = $wb.Worksheets.Item("2021-Q4")
$quarter2021.Range("A1:H3").Copy()
$newSheet.Range("A1:H3").PasteSpecial(-4122)

# Header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2. (B-G are text in the source data, incl. numeric-looking
# fund codes/percentages with significant leading zeros / trailing
# zeros, so they are entered with a leading quote - same as typing
# them into Excel directly - to keep them stored as text.)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'012864"
$newSheet.Range("C2").Value = "易方达标普医疗保健指数（QDII-LOF）人民币 C"
$newSheet.Range("D2").Value = "'0.50"
$newSheet.Range("E2").Value = "'93.65"
$newSheet.Range("F2").Value = "'1.53"
$newSheet.Range("G2").Value = "'0.0076"
$newSheet.Range("H2").Value = 10

# Row 3.
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'161126"
$newSheet.Range("C3").Value = "易方达标普医疗保健指数（QDII-LOF）人民币"
$newSheet.Range("D3").Value = "'0.50"
$newSheet.Range("E3").Value = "'93.65"
$newSheet.Range("F3").Value = "'1.53"
$newSheet.Range("G3").Value = "'0.0076"
$newSheet.Range("H3").Value = 10

# Row 4 (new - extends beyond the copied template range).
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'003719"
$newSheet.Range("C4").Value = "易方达标普医疗保健指数（QDII-LOF）美元A"
$newSheet.Range("D4").Value = "'0.49"
$newSheet.Range("E4").Value = "'93.65"
$newSheet.Range("F4").Value = "'1.53"
$newSheet.Range("G4").Value = "'0.0075"
$newSheet.Range("H4").Value = 10

# Row 5 (new - extends beyond the copied template range).
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'012865"
$newSheet.Range("C5").Value = "易方达标普医疗保健指数（QDII-LOF）美元 C"
$newSheet.Range("D5").Value = "'0.01"
$newSheet.Range("E5").Value = "'93.65"
$newSheet.Range("F5").Value = "'1.53"
$newSheet.Range("G5").Value = "'0.0002"
$newSheet.Range("H5").Value = 10

# Give rows 4 and 5 the same index-column / row styling as row 3 so the
# whole table looks consistent (A4/A5 pick up the centered border
# style used by A2/A3).
$newSheet.Range("A3").Copy()
$newSheet.Range("A4:A5").PasteSpecial(-4122)

Write-Output "done"
